# Fix typo in the Russian localization strings: replace the hyphenated
# "мини-игра"/"мини-игре" with the two-word form "мини игра"/"мини игре".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Replace("Перейти к мини-игре", "Перейти к мини игре") | Out-Null
$ws.Cells.Replace("Пройдена мини-игра", "Пройдена мини игра") | Out-Null
